$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, matching style of existing header H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Numeric data for columns I (I0) and J (IF), rows 2-7
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 11

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 6

$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 6

$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 5
